# Revert "nothing, just removed train files"
# This reverts a prior commit that blanked out the "done" markers in the
# results table (F7:I7, H8:I8, I10) and reset the view/selection.
# Restore the original "done" text in those cells and the selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 7: F7, G7, H7, I7 were cleared -> restore "done"
$ws.Range("F7").Value = "done"
$ws.Range("G7").Value = "done"
$ws.Range("H7").Value = "done"
$ws.Range("I7").Value = "done"

# Row 8: H8, I8 were cleared -> restore "done" (G8 already had it)
$ws.Range("H8").Value = "done"
$ws.Range("I8").Value = "done"

# Row 10: I10 was cleared -> restore "done"
$ws.Range("I10").Value = "done"

# Restore the sheet's prior scroll/selection state.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F17").Select()
